# Weekly fruit/vegetable price update:
# A new record (row) is inserted right after the existing row 509, shifting all
# subsequent rows down by one (509 stays, old 510 -> 511, ..., old 551 -> 552).
# Net effect: dimension grows from A1:R551 to A1:R552.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 510 (pushes existing rows 510..551 down to 511..552)
$ws.Rows.Item(510).Insert()

# Populate the newly inserted row 510 with the new price record
$ws.Cells.Item(510, 1).Value  = 10
$ws.Cells.Item(510, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(510, 3).Value  = "La Araucanía"
$ws.Cells.Item(510, 4).Value  = 44783
$ws.Cells.Item(510, 5).Value  = 9
$ws.Cells.Item(510, 6).Value  = 100112043
$ws.Cells.Item(510, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(510, 8).Value  = "Sin especificar"
$ws.Cells.Item(510, 9).Value  = "Primera"
$ws.Cells.Item(510, 10).Value = 50
$ws.Cells.Item(510, 11).Value = 21000
$ws.Cells.Item(510, 12).Value = 21000
$ws.Cells.Item(510, 13).Value = 21000
$ws.Cells.Item(510, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(510, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(510, 16).Value = 350
$ws.Cells.Item(510, 17).Value = 60
$ws.Cells.Item(510, 18).Value = "Hortaliza"
